# Updates the cryptos price/volume table to the latest scrape.
# Column D ("Price") values are written with a leading apostrophe so Excel
# keeps them as literal text (e.g. "1.000", "30.614.52") instead of
# re-parsing them as numbers and dropping the trailing zeros / thousands
# groups. Column E ("Volume(1h)") values already contain non-numeric
# characters (%, surrounding spaces) so Excel stores them as text as-is.
# A few rows (11-14, 44-45, 48-49) also changed rank order, so their
# Coin name (B) and Link (C) cells are rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.614.52"
$ws.Range("E2").Value = "  +0.51%  "

# Row 3
$ws.Range("D3").Value = "'1.882.09"
$ws.Range("E3").Value = "  +0.21%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'250.10"
$ws.Range("E5").Value = "  +1.07%  "

# Row 6
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").Value = "'0.4748"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("D8").Value = "'0.2933"
$ws.Range("E8").Value = "  +1.30%  "

# Row 9
$ws.Range("D9").Value = "'0.06528"
$ws.Range("E9").Value = "  +0.25%  "

# Row 10
$ws.Range("D10").Value = "'21.97"
$ws.Range("E10").Value = "  +1.83%  "

# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "'1.882.67"
$ws.Range("E11").Value = "  +0.26%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07736"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7403"
$ws.Range("E13").Value = "  -0.62%  "

# Row 14
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'96.82"
$ws.Range("E14").Value = "  -0.04%  "

# Row 15
$ws.Range("D15").Value = "'5.273"
$ws.Range("E15").Value = "  +3.07%  "

# Row 16
$ws.Range("D16").Value = "'274.57"
$ws.Range("E16").Value = "  +0.52%  "

# Row 17
$ws.Range("D17").Value = "'30.596.16"
$ws.Range("E17").Value = "  +0.45%  "

# Row 18
$ws.Range("D18").Value = "'13.19"
$ws.Range("E18").Value = "  -3.18%  "

# Row 19
$ws.Range("D19").Value = "'0.000007535"
$ws.Range("E19").Value = "  -0.17%  "

# Row 20
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.03%  "

# Row 21
$ws.Range("D21").Value = "'2.130.91"
$ws.Range("E21").Value = "  +0.35%  "

# Row 22
$ws.Range("D22").Value = "'5.332"
$ws.Range("E22").Value = "  +1.27%  "

# Row 23
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").Value = "'6.226"
$ws.Range("E24").Value = "  +1.08%  "

# Row 25
$ws.Range("D25").Value = "'9.213"
$ws.Range("E25").Value = "  -0.64%  "

# Row 26
$ws.Range("D26").Value = "'163.77"
$ws.Range("E26").Value = "  -0.36%  "

# Row 27
$ws.Range("D27").Value = "'18.86"
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("D28").Value = "'1.918"
$ws.Range("E28").Value = "  -1.74%  "

# Row 29
$ws.Range("E29").Value = "  -2.18%  "

# Row 30
$ws.Range("D30").Value = "'0.09691"
$ws.Range("E30").Value = "  -2.82%  "

# Row 31
$ws.Range("E31").Value = "  -0.47%  "

# Row 32
$ws.Range("D32").Value = "'4.292"
$ws.Range("E32").Value = "  -0.50%  "

# Row 33
$ws.Range("D33").Value = "'4.151"
$ws.Range("E33").Value = "  +2.36%  "

# Row 34
$ws.Range("D34").Value = "'0.04872"
$ws.Range("E34").Value = "  +2.23%  "

# Row 35
$ws.Range("D35").Value = "'1.128"
$ws.Range("E35").Value = "  +0.38%  "

# Row 36
$ws.Range("D36").Value = "'0.6990"
$ws.Range("E36").Value = "  +0.14%  "

# Row 37
$ws.Range("D37").Value = "'2.720"
$ws.Range("E37").Value = "  +0.17%  "

# Row 38
$ws.Range("D38").Value = "'0.01900"
$ws.Range("E38").Value = "  +1.76%  "

# Row 39
$ws.Range("D39").Value = "'2.770"
$ws.Range("E39").Value = "  +1.18%  "

# Row 40
$ws.Range("D40").Value = "'6.311"
$ws.Range("E40").Value = "  -0.74%  "

# Row 41
$ws.Range("D41").Value = "'74.85"
$ws.Range("E41").Value = "  +6.81%  "

# Row 42
$ws.Range("D42").Value = "'2.013"
$ws.Range("E42").Value = "  +4.64%  "

# Row 43
$ws.Range("D43").Value = "'0.4237"
$ws.Range("E43").Value = "  +1.62%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.8393"
$ws.Range("E45").Value = "  +0.28%  "

# Row 46
$ws.Range("D46").Value = "'102.65"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("D47").Value = "'9.369"
$ws.Range("E47").Value = "  +0.59%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.051"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'35.60"
$ws.Range("E49").Value = "  +0.81%  "

# Row 50
$ws.Range("D50").Value = "'918.02"
$ws.Range("E50").Value = "  -1.05%  "

# Row 51
$ws.Range("D51").Value = "'0.05732"
$ws.Range("E51").Value = "  +2.10%  "
